$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-16 (senders x targets with updated LR-pair statistics), columns A-T
$rowsData = @(
    @("ECs", "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 15.829186, 47.487558, 0.01520167221269649, 0.01552195334947967, 3, 1, 2.311298, 6.933894, 0.8122763614007964, 0.8638687058964239, 36.58596594342801, 329.273693490852, 0.0123479589921367, 0.01340892975299967),
    @("ECs", "Thbs1", "Tnfrsf11b", "M1", 3, 1, 15.829186, 47.487558, 0.01520167221269649, 0.01552195334947967, 1, 0.3333333333333333, 0.02434666666666667, 0.07303999999999999, 0.008556327142686946, 0.009099788701510982, 0.3853879151466666, 3.46849123632, 0.0001300704805677249, 0.0001412464957149756),
    @("ECs", "Thbs1", "Tnfrsf11b", "sCs", 3, 1, 15.829186, 47.487558, 0.01520167221269649, 0.01552195334947967, 2, 1, 0.509813, 1.019626, 0.1791673114565167, 0.1270315054020651, 8.069924802217999, 48.419548813308, 0.002723642739992067, 0.00197177710076503),
    @("FAPs", "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 155.500389, 466.501167, 0.1493359129474374, 0.1524822428572096, 3, 1, 2.311298, 6.933894, 0.8122763614007964, 0.8638687058964239, 359.4077380949221, 3234.669642854298, 0.1213020319954106, 0.1317246378092419),
    @("FAPs", "Thbs1", "Tnfrsf11b", "M1", 3, 1, 155.500389, 466.501167, 0.1493359129474374, 0.1524822428572096, 1, 0.3333333333333333, 0.02434666666666667, 0.07303999999999999, 0.008556327142686946, 0.009099788701510982, 3.78591613752, 34.07324523768, 0.001277766925330094, 0.001387556190733089),
    @("FAPs", "Thbs1", "Tnfrsf11b", "sCs", 3, 1, 155.500389, 466.501167, 0.1493359129474374, 0.1524822428572096, 2, 1, 0.509813, 1.019626, 0.1791673114565167, 0.1270315054020651, 79.276119817257, 475.656718903542, 0.02675611402669678, 0.01937004885723463),
    @("M1", "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 179.1193723333333, 537.358117, 0.1720185728536685, 0.1756427994052303, 3, 1, 2.311298, 6.933894, 0.8122763614007964, 0.8638687058964239, 413.9982470352887, 3725.984223317598, 0.1397266204509357, 0.1517323178222215),
    @("M1", "Thbs1", "Tnfrsf11b", "M1", 3, 1, 179.1193723333333, 537.358117, 0.1720185728536685, 0.1756427994052303, 1, 0.3333333333333333, 0.02434666666666667, 0.07303999999999999, 0.008556327142686946, 0.009099788701510982, 4.360959651742222, 39.24863686568, 0.001471847183954116, 0.001598312361529474),
    @("M1", "Thbs1", "Tnfrsf11b", "sCs", 3, 1, 179.1193723333333, 537.358117, 0.1720185728536685, 0.1756427994052303, 2, 1, 0.509813, 1.019626, 0.1791673114565167, 0.1270315054020651, 91.31738456737367, 547.904307404242, 0.03082010521877874, 0.02231216922147936),
    @("M2", "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 626.3728126666666, 1879.118438, 0.6015416194555684, 0.6142153852759307, 3, 1, 2.311298, 6.933894, 0.8122763614007964, 0.8638687058964239, 1447.734229170841, 13029.60806253757, 0.4886180378825116, 0.5306014500199917),
    @("M2", "Thbs1", "Tnfrsf11b", "M1", 3, 1, 626.3728126666666, 1879.118438, 0.6015416194555684, 0.6142153852759307, 1, 0.3333333333333333, 0.02434666666666667, 0.07303999999999999, 0.008556327142686946, 0.009099788701510982, 15.25009007905778, 137.25081071152, 0.005146986886003541, 0.005589230223228129),
    @("M2", "Thbs1", "Tnfrsf11b", "sCs", 3, 1, 626.3728126666666, 1879.118438, 0.6015416194555684, 0.6142153852759307, 2, 1, 0.509813, 1.019626, 0.1791673114565167, 0.1270315054020651, 319.3330027440313, 1915.998016464188, 0.1077765946870533, 0.07802470503271092),
    @("sCs", "Thbs1", "Tnfrsf11b", "FAPs", 2, 1, 64.45750049999999, 128.915001, 0.06190222253062919, 0.04213761911214986, 3, 1, 2.311298, 6.933894, 0.8122763614007964, 0.8638687058964239, 148.980491990649, 893.882951943894, 0.05028171207980187, 0.03640137049196932),
    @("sCs", "Thbs1", "Tnfrsf11b", "M1", 2, 1, 64.45750049999999, 128.915001, 0.06190222253062919, 0.04213761911214986, 1, 0.3333333333333333, 0.02434666666666667, 0.07303999999999999, 0.008556327142686946, 0.009099788701510982, 1.56932527884, 9.415951673039999, 0.0005296556668314699, 0.0003834434303053145),
    @("sCs", "Thbs1", "Tnfrsf11b", "sCs", 2, 1, 64.45750049999999, 128.915001, 0.06190222253062919, 0.04213761911214986, 2, 1, 0.509813, 1.019626, 0.1791673114565167, 0.1270315054020651, 32.8612717024065, 131.445086809626, 0.01109085478399584, 0.005352805189875229)
)

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $i + 2
    $rowVals = $rowsData[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}
